# Text updates as supplied by PM&C.
# Updates the "Notes" paragraph wording and splits the old combined
# "Sources: ..." sentence into a labelled "Source" row plus two separate
# lines (Census data / Survey data) on the Description sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# --- Row 7: revise the Notes text (comma added, "our" -> "the") ---
$ws.Range("B7").Value = "Assessment of progress to achieve this target uses Census data which is only available every five years. Trajectories for this target are not in a straight line from the baseline to target year. This is to reflect the lag between policy interventions and their anticipated impact. In the interim, survey data is used as a supplementary indicator, which is updated every three years. The small sample size for this age cohort has affected the ability to detect significant changes at the jurisdictional level."
$ws.Range("B7").Style = "Normal"

# --- Row 8: add a "Source" label in column A, and the Census data line in column B ---
$ws.Range("A8").Value = "Source"
$ws.Range("B8").Value = "Census data: ABS Census of Population and Housing 2006 and 2011 (Indigenous and non-Indigenous data)."
$ws.Range("B8").Style = "Normal"

# --- Row 9 (new row): the Survey data line in column B ---
$ws.Range("B9").Value = "Survey data: ABS National Aboriginal and Torres Strait Islander Social Survey (NATSISS) (2008 and 2014-15) and ABS Survey of Education and Work (SEW) (2008 and 2014)"
$ws.Range("B9").Style = "Normal"

# --- Row heights: rows 5/6 shrink slightly, rows 7-9 become single-line rows ---
$ws.Rows.Item(5).RowHeight = 73.45
$ws.Rows.Item(6).RowHeight = 109.45
$ws.Rows.Item(7).RowHeight = 13.8
$ws.Rows.Item(8).RowHeight = 13.8
$ws.Rows.Item(9).RowHeight = 13.8

# --- Update the selection to reflect the new last-used cell ---
$null = $ws.Range("B10").Select()
